$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2681.6667
$ws.Range("I64").Value = 2475
$ws.Range("J64").Value = 2888.3333
$ws.Range("K64").Value = 2475
$ws.Range("L64").Value = 2888.3333
$ws.Range("M64").Value = -2227
$ws.Range("N64").Value = -3384.3333
$ws.Range("H67").Value = 2681.6667
$ws.Range("I67").Value = 2475
$ws.Range("J67").Value = 2888.3333
$ws.Range("K67").Value = 2475
$ws.Range("L67").Value = 2888.3333
$ws.Range("M67").Value = -1617
$ws.Range("N67").Value = -4604.3333
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("H76").Value = 4885.143
$ws.Range("I76").Value = 4366
$ws.Range("K76").Value = 4366
$ws.Range("M76").Value = -4051
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0
$ws.Range("H79").Value = 4885.143
$ws.Range("I79").Value = 4366
$ws.Range("K79").Value = 4366
$ws.Range("M79").Value = -3274
$ws.Range("H92").Value = 207.5
$ws.Range("I92").Value = 207.5
$ws.Range("K92").Value = 207.5
$ws.Range("M92").Value = 1040.5
$ws.Range("H94").Value = 2220.889
$ws.Range("I94").Value = 1873.5
$ws.Range("K94").Value = 1873.5
$ws.Range("M94").Value = -1422.5
$ws.Range("H100").Value = 4471
$ws.Range("I100").Value = 2452
$ws.Range("J100").Value = 6490
$ws.Range("K100").Value = 2452
$ws.Range("L100").Value = 6490
$ws.Range("M100").Value = -1911
$ws.Range("N100").Value = -7572
$ws.Range("H112").Value = 1048.8125
$ws.Range("J112").Value = 1048.8125
$ws.Range("L112").Value = 3146.4375
$ws.Range("N112").Value = -5362.4375
$ws.Range("H129").Value = 922.81134
$ws.Range("J129").Value = 899.1277
$ws.Range("L129").Value = 2697.3831
$ws.Range("N129").Value = -12697.3831
$ws.Range("H135").Value = 604.8333
$ws.Range("I135").Value = 590.375
$ws.Range("J135").Value = 633.75
$ws.Range("K135").Value = 5313.375
$ws.Range("L135").Value = 5703.75
$ws.Range("M135").Value = -2778.375
$ws.Range("N135").Value = -10773.75
$ws.Range("H138").Value = 1895.5625
$ws.Range("J138").Value = 2500
$ws.Range("L138").Value = 7500
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5999.5
$ws.Range("I63").Value = 5999.5
$ws.Range("K63").Value = 5999.5
$ws.Range("M63").Value = -5313.5
$ws.Range("H66").Value = 5999.5
$ws.Range("I66").Value = 5999.5
$ws.Range("K66").Value = 29997.5
$ws.Range("M66").Value = -26565.5
$ws.Range("H74").Value = 1087.4814
$ws.Range("I74").Value = 482.4
$ws.Range("K74").Value = 482.4
$ws.Range("M74").Value = 391.6
$ws.Range("H77").Value = 1087.4814
$ws.Range("I77").Value = 482.4
$ws.Range("K77").Value = 2412
$ws.Range("M77").Value = 1956
$ws.Range("H97").Value = 1633.1875
$ws.Range("I97").Value = 1542.1333
$ws.Range("K97").Value = 1542.1333
$ws.Range("M97").Value = -1046.1333
$ws.Range("H102").Value = 3599.8
$ws.Range("I102").Value = 1999.75
$ws.Range("K102").Value = 1999.75
$ws.Range("M102").Value = -377.75
$ws.Range("H132").Value = 2193.3696
$ws.Range("I132").Value = 1616.0322
$ws.Range("K132").Value = 4848.096600000001
$ws.Range("M132").Value = -2318.096600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1750.28
$ws.Range("I20").Value = 1778.125
$ws.Range("J20").Value = 1700.7778
$ws.Range("K20").Value = 1778.125
$ws.Range("L20").Value = 1700.7778
$ws.Range("M20").Value = -1531.125
$ws.Range("N20").Value = -2194.7778
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").Value = 0
$ws.Range("H82").Value = 39499.75
$ws.Range("H85").Value = 39499.75
$ws.Range("H94").Value = 1624.6154
$ws.Range("J94").Value = 1690
$ws.Range("L94").Value = 1690
$ws.Range("N94").Value = -2592
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2
$ws.Range("H105").Value = 2452.45
$ws.Range("I105").Value = 2585.2354
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 2585.2354
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = -838.2354
$ws.Range("N105").Value = -5194
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1372.25
$ws.Range("J31").Value = 1564.85
$ws.Range("L31").Value = 1564.85
$ws.Range("N31").Value = -2154.85
$ws.Range("H34").Value = 1372.25
$ws.Range("J34").Value = 1564.85
$ws.Range("L34").Value = 1564.85
$ws.Range("N34").Value = -1968.85
$ws.Range("H62").Value = 2904
$ws.Range("I62").Value = 2904
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2904
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -2280
$ws.Range("H65").Value = 2904
$ws.Range("I65").Value = 2904
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14520
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -11400
$ws.Range("H68").Value = 43000
$ws.Range("J68").Value = 43000
$ws.Range("L68").Value = 43000
$ws.Range("N68").Value = -44498
$ws.Range("H71").Value = 43000
$ws.Range("J71").Value = 43000
$ws.Range("L71").Value = 129000
$ws.Range("N71").Value = -136488
$ws.Range("H132").Value = 3482.2
$ws.Range("I132").Value = 2847.8333
$ws.Range("J132").Value = 3905.111
$ws.Range("K132").Value = 8543.499899999999
$ws.Range("L132").Value = 11715.333
$ws.Range("M132").Value = -6013.499899999999
$ws.Range("N132").Value = -16775.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 4812
$ws.Range("I44").Value = 2050
$ws.Range("J44").Value = 6653.3335
$ws.Range("K44").Value = 6150
$ws.Range("L44").Value = 19960.0005
$ws.Range("M44").Value = -5752
$ws.Range("N44").Value = -20756.0005
$ws.Range("H68").Value = 1621.7273
$ws.Range("J68").Value = 1778.8379
$ws.Range("L68").Value = 5336.5137
$ws.Range("N68").Value = -6958.5137
$ws.Range("H71").Value = 1621.7273
$ws.Range("J71").Value = 1778.8379
$ws.Range("L71").Value = 16009.5411
$ws.Range("N71").Value = -24121.5411
$ws.Range("H121").Value = 616.25
$ws.Range("I121").Value = 432.5
$ws.Range("K121").Value = 1297.5
$ws.Range("M121").Value = 12.5
$ws.Range("H131").Value = 11380813
$ws.Range("J131").Value = 19838.658
$ws.Range("L131").Value = 59515.974
$ws.Range("N131").Value = -69595.974

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 20000.5
$ws.Range("J15").Value = 20000.5
$ws.Range("L15").Value = 20000.5
$ws.Range("N15").Value = -20576.5
$ws.Range("H81").Value = 20000.5
$ws.Range("J81").Value = 20000.5
$ws.Range("L81").Value = 20000.5
$ws.Range("N81").Value = -21996.5
$ws.Range("H84").Value = 20000.5
$ws.Range("J84").Value = 20000.5
$ws.Range("L84").Value = 60001.5
$ws.Range("N84").Value = -69985.5
$ws.Range("H102").Value = 3384.724
$ws.Range("I102").Value = 3339.3333
$ws.Range("J102").Value = 3602.6
$ws.Range("K102").Value = 3339.3333
$ws.Range("L102").Value = 3602.6
$ws.Range("M102").Value = -1717.3333
$ws.Range("N102").Value = -6846.6
$ws.Range("H132").Value = 5499349.5
$ws.Range("I132").Value = 19236274
$ws.Range("J132").Value = 4579.4
$ws.Range("K132").Value = 57708822
$ws.Range("L132").Value = 13738.2
$ws.Range("M132").Value = -57706292
$ws.Range("N132").Value = -18798.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2728.25
$ws.Range("J22").Value = 4496.5
$ws.Range("L22").Value = 4496.5
$ws.Range("N22").Value = -5086.5
$ws.Range("H27").Value = 2728.25
$ws.Range("J27").Value = 4496.5
$ws.Range("L27").Value = 4496.5
$ws.Range("N27").Value = -4710.5
$ws.Range("H40").Value = 7110.625
$ws.Range("I40").Value = 3263.6667
$ws.Range("K40").Value = 3263.6667
$ws.Range("M40").Value = -3127.6667
$ws.Range("H61").Value = 3800.3333
$ws.Range("I61").Value = 2303.75
$ws.Range("K61").Value = 2303.75
$ws.Range("M61").Value = -2101.75
$ws.Range("H68").Value = 1399.6
$ws.Range("I68").Value = 1499.5
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1499.5
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -750.5
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 1399.6
$ws.Range("I71").Value = 1499.5
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 7497.5
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -3753.5
$ws.Range("N71").Value = -12488
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H113").Value = 3800.3333
$ws.Range("I113").Value = 2303.75
$ws.Range("K113").Value = 2303.75
$ws.Range("M113").Value = -133.75
$ws.Range("H132").Value = 4737.7
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 4897.579
$ws.Range("K132").Value = 5100
$ws.Range("L132").Value = 14692.737
$ws.Range("M132").Value = -2570
$ws.Range("N132").Value = -19752.737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -2280
$ws.Range("H123").Value = 41549.332
$ws.Range("J123").Value = 41549.332
$ws.Range("L123").Value = 41549.332
$ws.Range("N123").Value = -51349.332
$ws.Range("H124").Value = 17499.5
$ws.Range("J124").Value = 17499.5
$ws.Range("L124").Value = 17499.5
$ws.Range("N124").Value = -27319.5
$ws.Range("H132").Value = 2349.0454
$ws.Range("I132").Value = 1781.091
$ws.Range("J132").Value = 2917
$ws.Range("K132").Value = 5343.272999999999
$ws.Range("L132").Value = 8751
$ws.Range("M132").Value = -2813.272999999999
$ws.Range("N132").Value = -13811
